# Highlight key phrases (yellow) in Task 2, Task 3 and Task 4 paragraphs.
$d = $word.ActiveDocument

$phrases = @(
    "seamless cloning",
    "importing gradients",
    "mixing gradients",
    "2a",
    "only one",
    "texture flattening",
    "local illumination changes",
    "local colour changes",
    "seamless tiling"
)

foreach ($phrase in $phrases) {
    $rng = $d.Content
    $found = $rng.Find.Execute($phrase, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Font.HighlightColorIndex = 7
    }
}
